# Apply the weekly cryptos price/volume refresh (GitHub Actions data update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.293.12"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.78%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.712.20"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.71%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.42%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.55"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.37%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5267"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.59%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.007"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.35%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06653"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.77%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2640"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.82%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.72"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.86%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07736"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.70%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.478"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.78%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.722.58"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.15%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.949.64"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.63%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5772"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.60%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8160"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.43%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.61"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.40%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.327.36"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.64%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.78"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.29%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.32%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.638"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.08%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.39"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.73%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.016"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.03%  "

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.33%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.41"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.26%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.707"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.31%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1203"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.37%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.217"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.48%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.15"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.12%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05367"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.68%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.294"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.35%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.470"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.27%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.370"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.93%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.630"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.32%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.849"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.46%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9494"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.97%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.397"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.13%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5864"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.20%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.153.48"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +10.09%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01648"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.30%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.829"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.30%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.30%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8381"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.85%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.91"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.31%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.856.56"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.65%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.79%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.40"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.21%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4557"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.17%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.005"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.18%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.098"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.11%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05235"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.28%  "
